$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Fill in the previously-empty "Description" cell for the Multiplan URL row.
$ws.Range("C6").Value = "Multiplan Website URL"

# 2) Insert a brand-new row right after the Multiplan row (row 7) for the
#    Delta Dental website URL setting. This pushes MaxRetry and the
#    Aetna..Delta Dental table down by one row (old row 7 -> 8, ... old row 18 -> 19).
$ws.Rows("7:7").Insert()
$ws.Rows("7:7").RowHeight = 14.25

$ws.Range("A7").Value = "URL_Delta"
$ws.Range("B7").Value = "https://www.deltadentalins.com/find-a-dentist/directory-results/?d=enterprise&network=2ppo&location=2064%20Baldwin%20St;%20Ste%20A;%20Jenison;%20MI%2049428&distance=15&lat=42.9068881&long=-85.8330607&free_text=Vurugonda%20Anupama&sort_field=relevance&sort_order=asc&isBookmarkedURL=false&page=1"
$ws.Range("C7").Value = "Delta Dental Website"

# 3) Flip the PPO and PHCS rows (now rows 15 & 16) from "N" to "Y", and the
#    Anthem row (now row 12) from "Y" back to "N" - matches shared-string
#    index swap seen in the diff (B12 <v>113</v> = "N", B15/B16 <v>112</v> = "Y").
$ws.Range("B12").Value = "N"
$ws.Range("B15").Value = "Y"
$ws.Range("B16").Value = "Y"

# 4) The Delta Dental row (now row 19, previously blank in column B) gets
#    marked "Y" and the row is fully populated (it already carries the
#    Name/Description from the row shift, so just set the Value column).
$ws.Range("B19").Value = "Y"

# 5) The conditional formatting range that highlights the Y/N table shifts
#    down by one row along with the table (A8:A18 -> A9:A19), and the rule
#    formulas move from $B8 to $B9. Reuse the existing rules (same dxfId /
#    priority) instead of recreating them.
$fcs = $ws.Range("A8:A18").FormatConditions
$fcN = $fcs.Item(1)
$fcY = $fcs.Item(2)
$fcN.ModifyAppliesToRange($ws.Range("A9:A19"))
$fcN.Formula1 = '=$B9="N"'
$fcY.Formula1 = '=$B9="Y"'

# 6) Restore the active-cell selection to B7, matching the authored workbook.
$ws.Range("B7").Select()
